# Update cryptocurrency price values in column D to reflect the latest
# symbol-list data refresh (commit: "Updated symbol list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "267.85"
    "D3"  = "22.86"
    "D4"  = "6.326"
    "D5"  = "0.06202"
    "D6"  = "3.589"
    "D7"  = "6.678"
    "D8"  = "1.400"
    "D9"  = "0.8319"
    "D12" = "0.08225"
    "D13" = "0.03400"
    "D14" = "0.03153"
    "D15" = "0.09279"
    "D16" = "3.912"
    "D17" = "0.001726"
    "D18" = "0.04841"
    "D19" = "0.006300"
    "D20" = "0.005380"
    "D22" = "0.0001500"
    "D24" = "2.369"
    "D25" = "0.3350"
    "D26" = "0.1212"
    "D40" = "0.04652"
    "D41" = "0.006868"
    "D42" = "0.1154"
    "D44" = "0.01213"
    "D45" = "0.00006254"
    "D47" = "0.7003"
    "D48" = "0.1658"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
